$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Result" values for the UL (row 2) and NASR (row 3) task rows.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "N"

# Move the active selection from C5 to G15.
$ws.Range("G15").Select()
